$d = $word.ActiveDocument

# 1) Simple text replacements (paragraph content swaps; Word will add/drop
#    xml:space="preserve" automatically based on leading/trailing spaces).
$d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק 08.09.24: ⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק 07.09.24: ⚡️🚀", 2) | Out-Null
$d.Content.Find.Execute("DO TRANSFORMER WORLD MODELS GIVE BETTER POLICY GRADIENTS?", $true, $false, $false, $false, $false, $true, 1, $false, "ReMamba: Equip Mamba with Effective Long-Sequence Modeling", 2) | Out-Null
$d.Content.Find.Execute("לא הייתי אמור לכתוב סקירה היום אך הקלטת הפודקאסט שלנו התבטלה והתפנה לי קצת זמן אז אסקור מאמר שכבר נמצא כמה זמן אצלי במגירה. המאמר בנושא למידה עם חיזוקים (RL) וטרנספורמרים אז לכאורה זה נשמע מאמר די נחמד.", $true, $false, $false, $false, $false, $true, 1, $false, " סוקר את המאמר הזה משתי סיבות. קודם כל הוא קשור לממבה. הסיבה השני היא זה שהתבקשתי לסקור אותו. ואוקיי, המאמר לצערי לא חידש לי הרבה ולדעתי לא נזכור אותו בעוד כמה חודשים. ", 2) | Out-Null
$d.Content.Find.Execute("המאמר מדבר על שיטה לשיפור של למידה פוליסי בבעיות של RL בבעיות שיש לנו גישה ישירה לדינמיקה של הסביבה (כלומר אנו לא יכולים לאסוף עליה דאטה רלוונטי המאפיין את פיצ'רים המהותיים שלו קרי non-observable). בגדול המטרה שלנו בלמידת פוליסי היא לחזות את הפעולה (action) האופטימלי בהינתן המצב s של הסביבה והפעולה האחרונה s. אופטימלי כאן משמעותו מקסום של התגמול (reward) הכולל המתקבל במהלך אפיזודה. המודל שחוזה את הפעולה הזו הוא למעשה מממש את הפוליסי שלנו. ", $true, $false, $false, $false, $false, $true, 1, $false, "אתם זוכרים את State Space Models או SSM בהקשר של למודלים עמוקים? SSM הוא ארכיטקטורה יחסית חדשה עבור רשתות לעיבוד דאטה סדרתי (שפה טבעית וגם תמונות). השוס הגדול ב-SSM היא שהם מאוד מהירים גם באימון וגם באינפרנס עקב כך שניתן לייצג אותם בתור רשת קונבולוציה וגם במודל רשת recurrent. הגמישות הזו כמובן גובה מאתנו מחיר בדמות חוסר expressiveness (יכולת למדל חוקיות מורכבות) של ארכיטקטורה הזו עקב העובדה המעברים בין המצבים החבויים הם לינאריים וקבועים לכל איברי הסדרה (מכאן בא הדואליות בייצוג).", 2) | Out-Null
$d.Content.Find.Execute("אבל מה לעשות אם אין לנו גישה ישירה לסביבה? במקרה הזה אנו יכולים לאמן מודל שהוא חוזה לנו את המצב הבא s בהינתן המצב הקודם(כלומר ייצוגו) והפעולה האחרונה(עם הנחת המרקוביות) או בהינתן N ייצוגים של המצבים האחרונים והפעולה האחרונה. זה למעשה נקרא world model (לדעתי יחד עם מודלים המשערכים את התגמול הצפוי למצב נתון -   value function אבל זה פחות חשוב כרגע). ", $true, $false, $false, $false, $false, $true, 1, $false, "ארכיטקטורת ממבה מחזירה לנו קצת מה-expressiveness בכך שהופכת את המעברים בין המצבים החבויים לתלוי במצב החבוי אך משאיר אותם לינאריים. זה עוזר אבל עדיין ממבה מתקשה במשימות reasoning מורכבות עקב מחסור ב-expressiveness. ייתכן שאחת הסיבות לאי הצלחה זו היא חוסר יכולת של ארכיטקטורת ממבה לדחוס את המידע הרלוונטי למשימה (לגיטימי אבל כמובן יש עוד סיבות לכך).", 2) | Out-Null
$d.Content.Find.Execute("איך המודל הזה מאומן? מאינטראקציה עם הסביבה - הסוכן מבצע פעולות בסביבה ואנו מעדכנים את ה-world model שלנו בהתבסס על משוואות Bellman). שימו לב אם או ללא הנחת מקרוביות אנחנו משערכים את הייצוג של המצב ה״עולם״ הבא בהינתן המצב(-ים) הקודמים. המאמר טוען שזה יוצר גרדיאנטים לא יציבים ושונות גבוהה עקב שימוש ישיר בשערוך של המצבים הקודמים לשעורך של המצב הבא.", $true, $false, $false, $false, $false, $true, 1, $false, "המחברים מציעים לדחוס את ייצוגיהם של תת סדרות של טוקנים. נניח שיש לנו L טוקנים בהקשר ואנו רוצים ״לדחוס״ את טוקנים שייצוגיהם דומה לזה של הטוקן L. כלומר מחשבים את הדמיון בין תת-סדרה רציפה נתונה של טוקנים (הייפרפרמטר) ודוחסים את הייצוגים של הטוקנים בתת-סדרה זו לפחות טוקנים (הייפרפרמטר גם כן). כלומר במקום הייצוגים של M טוקנים בתת סדרה נקבל ייצוגים של K טוקנים אחרי הדחיסה.הדמיון מחושב דרך דמיון קוסיין (עם כל מיני שכבות לינאריות).", 2) | Out-Null
$d.Content.Find.Execute("הם מציע לשערך את המצב הבא מהפעולה ולא מייצוגי המצבים שטענתם ״הופך את הגרדיאנטים במודל העולם לפחות מעגליים״ וזה תורם ליציבות השערוך. יש גם קצת הוכחות במאמר (סוג של) של הטענה הזו. המאמר מראה אם יש לנו מקרוביות (התלות של המצב הבא היא רק במצב האחרון) השיטה המוצעת עובדת כמו RNN מבחינת הגרדיאנטים. בתחושה זה נשמע לי די טבעי (אשמח אם מישהו ירחיב על זה). במקרה שאין לנו מרקוביות הטענה לביצועים טובים יותר של השיטה המוצעת. ", $true, $false, $false, $false, $false, $true, 1, $false, "הם מראים שזה עובד - לי זה מריח קצת אוברפיט וגם קושי באופטימיזציה של ההייפרפרמטרים….", 2) | Out-Null

# 2) Delete the paragraph that only said the transformers were barely mentioned
#    (the whole <w:p> must go, not just its text).
$deleteTarget = "לא ראיתי אזכור משמעותי מדי של הטרנספורמרים במאמר (תקנו אותי אם אני טועה)."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $deleteTarget) {
        $p.Range.Delete()
        break
    }
}

# 3) Swap the old arxiv link for the new one (on the now-last paragraph).
$d.Content.Find.Execute("https://arxiv.org/abs/2402.05290", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2408.15496", 2) | Out-Null
